$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated raw counts for row 8 (ano = 2025)
$ws.Range("C8").Value = 1047
$ws.Range("D8").Value = 174
$ws.Range("E8").Value = 873

# Recalculated metrics
# F8 = returning_customers (D8) / total_customers of previous row (C7) * 100
$c7 = $ws.Range("C7").Value()
$d8 = $ws.Range("D8").Value()
$c8 = $ws.Range("C8").Value()
$e8 = $ws.Range("E8").Value()

$ws.Range("F8").Value = ($d8 / $c7) * 100
$ws.Range("G8").Value = ($e8 / $c8) * 100
$ws.Range("H8").Value = ($d8 / $c8) * 100
